$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 17
$ws.Range("H17").Value = 1266.712068611907
$ws.Range("I17").Value = 1414.634209378984

# Row 18
$ws.Range("H18").Value = 1211.184666308004
$ws.Range("I18").Value = 1225.698890613596

# Row 21
$ws.Range("H21").Value = 826.0106721316131
$ws.Range("I21").Value = 949.1475679611532

# Row 22
$ws.Range("H22").Value = 1714.950834610086
$ws.Range("I22").Value = 2425.306729108642
$ws.Range("J22").Value = 0

# Row 24
$ws.Range("H24").Value = 1251.886692580295
$ws.Range("I24").Value = 3208.981321493234
